# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# The source data pipeline recomputed the "K" (strikeouts) column (column G)
# for every saved game row from the underlying play-by-play log instead of the
# previous "Strike#" derived figure. This updates the already-written G2:G40
# values on the active sheet to the freshly regenerated strikeout counts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @(2, 0, 0, 2, 1, 1, 1, 1, 1, 1, 1, 1, 0, 0, 2, 1, 1, 1, 0, 0, 0, 4, 1, 0, 1, 0, 0, 0, 1, 0, 2, 0, 1, 0, 1, 0, 1, 1, 1)

$firstRow = 2
for ($i = 0; $i -lt $newK.Length; $i++) {
    $row = $firstRow + $i
    $ws.Cells.Item($row, 7).Value = $newK[$i]
}
